# Historias de usuario Tercer Sprint
#
# Applies, in order:
#  1. Appends a stray run "   3 3 3 3  33" to the
#     "Registro Materia Prima (Cantidad / Costo)" paragraph.
#  2. Appends a stray run " 3 3  3 3 3 3 " to the "Registro Mercancía"
#     paragraph, and moves the "_GoBack" bookmark (which used to sit at
#     the end of the "Registrar Proveedores" paragraph) onto the end of
#     this paragraph instead.
#  3. Removes the (now relocated) "_GoBack" bookmark from the end of the
#     "Registrar Proveedores" paragraph.
#  4. Wraps the "Priodidad" run in a spellcheck proofErr pair
#     (spellStart/spellEnd) - Word flagging it as a misspelling.
#  5. Splits the "*Subsitema de control de inventario " run into three
#     runs so the misspelled word "Subsitema" is wrapped in its own
#     spellStart/spellEnd proofErr pair.

$d = $word.ActiveDocument

function Insert-ParagraphXml($paragraph, $innerXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>'
    $pkg += '<?mso-application progid="Word.Document"?>'
    $pkg += '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
    $pkg += '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">'
    $pkg += '<pkg:xmlData>'
    $pkg += '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
    $pkg += '<w:body>'
    $pkg += $innerXml
    $pkg += '<w:sectPr/></w:body></w:document>'
    $pkg += '</pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($pkg)
}

# 1) "Registro Materia Prima (Cantidad / Costo)" -> add trailing run.
$pMateriaPrima = $d.Paragraphs(7)
$xmlMateriaPrima = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$xmlMateriaPrima += '<w:r><w:t>Registro Materia Prima (Cantidad / Costo)</w:t></w:r>'
$xmlMateriaPrima += '<w:r><w:t xml:space="preserve">   3 3 3 3  33</w:t></w:r>'
$xmlMateriaPrima += '</w:p>'
Insert-ParagraphXml $pMateriaPrima $xmlMateriaPrima

# 2) "Registro Mercancía" -> add trailing run + relocate the _GoBack bookmark here.
$pMercancia = $d.Paragraphs(8)
$xmlMercancia = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$xmlMercancia += '<w:r><w:t>Registro Mercancía</w:t></w:r>'
$xmlMercancia += '<w:r><w:t xml:space="preserve"> 3 3  3 3 3 3 </w:t></w:r>'
$xmlMercancia += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$xmlMercancia += '</w:p>'
Insert-ParagraphXml $pMercancia $xmlMercancia

# 3) "Registrar Proveedores" -> drop the bookmark that used to close this paragraph.
$pRegistrarProveedores = $d.Paragraphs(28)
$xmlRegistrarProveedores = '<w:p><w:r><w:t>Registrar Proveedores</w:t></w:r></w:p>'
Insert-ParagraphXml $pRegistrarProveedores $xmlRegistrarProveedores

# 4) "Priodidad" -> wrap with proofErr spellStart/spellEnd.
$pPriodidad = $d.Paragraphs(24)
$xmlPriodidad = '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Priodidad</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Insert-ParagraphXml $pPriodidad $xmlPriodidad

# 5) "*Subsitema de control de inventario " -> split into 3 runs, proofErr around "Subsitema".
$pSubsistema = $d.Paragraphs(32)
$xmlSubsistema = '<w:p><w:pPr><w:tabs><w:tab w:val="center" w:pos="4419"/></w:tabs></w:pPr>'
$xmlSubsistema += '<w:r><w:t>*</w:t></w:r>'
$xmlSubsistema += '<w:proofErr w:type="spellStart"/><w:r><w:t>Subsitema</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$xmlSubsistema += '<w:r><w:t xml:space="preserve"> de control de inventario </w:t></w:r>'
$xmlSubsistema += '<w:r><w:tab/></w:r>'
$xmlSubsistema += '</w:p>'
Insert-ParagraphXml $pSubsistema $xmlSubsistema
